$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Gas any)
$ws.Range("C2").Value = 1077
$ws.Range("D2").Value = 906
$ws.Range("E2").Value = 5169
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0.1724303554274736
$ws.Range("H2").Value = 0.1632642094607124
$ws.Range("I2").Value = 0.181999182105042
$ws.Range("J2").Value = 0.5431164901664145
$ws.Range("K2").Value = 0.5211293568207866
$ws.Range("L2").Value = 0.5649368963474446

# Row 3 (ABG threshold)
$ws.Range("C3").Value = 526
$ws.Range("D3").Value = 1457
$ws.Range("E3").Value = 2986
$ws.Range("F3").Value = 2183
$ws.Range("G3").Value = 0.1497722095671982
$ws.Range("H3").Value = 0.13835315820706
$ws.Range("I3").Value = 0.1619565887446273
$ws.Range("J3").Value = 0.2652546646495209
$ws.Range("K3").Value = 0.246291400556881
$ws.Range("L3").Value = 0.2851256655271927

# Row 4 (VBG threshold)
$ws.Range("C4").Value = 865
$ws.Range("D4").Value = 1118
$ws.Range("E4").Value = 2642
$ws.Range("F4").Value = 2527
$ws.Range("G4").Value = 0.2466495580268035
$ws.Range("H4").Value = 0.2326653313210056
$ws.Range("I4").Value = 0.2611882018439748
$ws.Range("J4").Value = 0.4362077660110943
$ws.Range("K4").Value = 0.4145249211005427
$ws.Range("L4").Value = 0.4581372891242377

# Row 5 - rename labels and update values
$ws.Range("A5").Value = "PCO2 UNKNOWN threshold"
$ws.Range("B5").Value = "unknown_hypercap_threshold"
$ws.Range("C5").Value = 69
$ws.Range("D5").Value = 1914
$ws.Range("E5").Value = 327
$ws.Range("F5").Value = 4842
$ws.Range("G5").Value = 0.1742424242424243
$ws.Range("H5").Value = 0.1400608016450688
$ws.Range("I5").Value = 0.2146834493364229
$ws.Range("J5").Value = 0.03479576399394856
$ws.Range("K5").Value = 0.0275869478302369
$ws.Range("L5").Value = 0.04380347851252991
